# "Final Test Case data"
#
# The sheet originally had 16 data columns (A:P) per row. The edit trims
# the table down to columns A:H, pulling in the values that used to live
# in columns K:N (so the new F/G/H values are the old K/L/M/N values),
# dropping the old columns in between, and - for the last (styled) data
# row only - leaving behind 8 empty-but-styled placeholder cells in I:P
# instead of removing them outright.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete columns E:J outright (6 columns). This shifts the old K:P
#    block left by six columns, so old K/L/M/N (the "R21..R24" data)
#    become the new E/F/G/H, exactly matching the target values.
$ws.Range("E1:J5").EntireColumn.Delete()

# 2) Rows 1-4 should end up with only columns A:H populated - clear what
#    is now I:J (the old O:P / "R25"+"R26" columns) for those rows so the
#    cells disappear entirely (they carry no style, so ClearContents drops
#    the <c> elements completely).
$ws.Range("I1:J4").ClearContents()

# 3) Row 5 keeps placeholders all the way out to column P (matching the
#    original dimension), but with no values - just the row's number
#    style carried along. Copy the already-styled I5:J5 cells' formatting
#    out to K5, M5 and O5 (which also stamps their neighbour L5/N5/P5),
#    recreating styled-but-empty cells across I5:P5, then clear all of
#    their contents.
$ws.Range("I5:J5").Copy()
$ws.Range("K5").PasteSpecial(-4122)
$ws.Range("M5").PasteSpecial(-4122)
$ws.Range("O5").PasteSpecial(-4122)
$ws.Range("I5:P5").ClearContents()

# 4) Update the active selection to match the saved view.
$selectionResult = $ws.Range("M13").Select()
